{"js": "// Update the date line and the 25 multiplication problems in the practice\n// table to the new values from the latest generated output.\n\nconst replacements = [\n  [\"2025-02-05 Wednesday\", \"2025-02-06 Thursday\"],\n  [\"125\u00d75=\", \"960\u00d77=\"],\n  [\"760\u00d78=\", \"765\u00d75=\"],\n  [\"788\u00d78=\", \"409\u00d74=\"],\n  [\"888\u00d77=\", \"364\u00d72=\"],\n  [\"102\u00d77=\", \"817\u00d72=\"],\n  [\"874\u00d76=\", \"574\u00d76=\"],\n  [\"887\u00d73=\", \"958\u00d79=\"],\n  [\"770\u00d72=\", \"698\u00d72=\"],\n  [\"925\u00d76=\", \"579\u00d72=\"],\n  [\"300\u00d77=\", \"985\u00d78=\"],\n  [\"759\u00d74=\", \"936\u00d75=\"],\n  [\"718\u00d75=\", \"445\u00d74=\"],\n  [\"467\u00d72=\", \"995\u00d75=\"],\n  [\"711\u00d72=\", \"491\u00d73=\"],\n  [\"738\u00d79=\", \"176\u00d79=\"],\n  [\"127\u00d79=\", \"440\u00d72=\"],\n  [\"317\u00d77=\", \"807\u00d73=\"],\n  [\"779\u00d73=\", \"304\u00d76=\"],\n  [\"265\u00d79=\", \"111\u00d77=\"],\n  [\"333\u00d75=\", \"670\u00d78=\"],\n  [\"985\u00d77=\", \"473\u00d79=\"],\n  [\"359\u00d75=\", \"589\u00d77=\"],\n  [\"714\u00d76=\", \"589\u00d76=\"],\n  [\"318\u00d78=\", \"494\u00d73=\"],\n  [\"404\u00d77=\", \"711\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in the practice\n# table to the new values from the latest generated output.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-02-05 Wednesday\"; New = \"2025-02-06 Thursday\" },\n    @{ Old = \"125\u00d75=\"; New = \"960\u00d77=\" },\n    @{ Old = \"760\u00d78=\"; New = \"765\u00d75=\" },\n    @{ Old = \"788\u00d78=\"; New = \"409\u00d74=\" },\n    @{ Old = \"888\u00d77=\"; New = \"364\u00d72=\" },\n    @{ Old = \"102\u00d77=\"; New = \"817\u00d72=\" },\n    @{ Old = \"874\u00d76=\"; New = \"574\u00d76=\" },\n    @{ Old = \"887\u00d73=\"; New = \"958\u00d79=\" },\n    @{ Old = \"770\u00d72=\"; New = \"698\u00d72=\" },\n    @{ Old = \"925\u00d76=\"; New = \"579\u00d72=\" },\n    @{ Old = \"300\u00d77=\"; New = \"985\u00d78=\" },\n    @{ Old = \"759\u00d74=\"; New = \"936\u00d75=\" },\n    @{ Old = \"718\u00d75=\"; New = \"445\u00d74=\" },\n    @{ Old = \"467\u00d72=\"; New = \"995\u00d75=\" },\n    @{ Old = \"711\u00d72=\"; New = \"491\u00d73=\" },\n    @{ Old = \"738\u00d79=\"; New = \"176\u00d79=\" },\n    @{ Old = \"127\u00d79=\"; New = \"440\u00d72=\" },\n    @{ Old = \"317\u00d77=\"; New = \"807\u00d73=\" },\n    @{ Old = \"779\u00d73=\"; New = \"304\u00d76=\" },\n    @{ Old = \"265\u00d79=\"; New = \"111\u00d77=\" },\n    @{ Old = \"333\u00d75=\"; New = \"670\u00d78=\" },\n    @{ Old = \"985\u00d77=\"; New = \"473\u00d79=\" },\n    @{ Old = \"359\u00d75=\"; New = \"589\u00d77=\" },\n    @{ Old = \"714\u00d76=\"; New = \"589\u00d76=\" },\n    @{ Old = \"318\u00d78=\"; New = \"494\u00d73=\" },\n    @{ Old = \"404\u00d77=\"; New = \"711\u00d73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
